# Applies the edits described by the commit diff:
#  - Rename sheet "Calls-frequency" -> "Call-frequency"
#  - Move the active selection on that sheet from G17 to C21

$wb = $excel.ActiveWorkbook

# The sheet that needs renaming is the one bound to rId1 / sheetId 2,
# which is the first sheet in the workbook and currently tab-selected.
$ws = $wb.Worksheets.Item("Calls-frequency")

# Rename the worksheet.
$ws.Name = "Call-frequency"

# Make sure it is the active sheet, then move the selection to C21.
$ws.Activate()
$ws.Range("C21").Select()
